$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These three products were removed from the stock/sales report.
$targets = @("MEPAFURAN 100 MG 20 CAPS.", "PERLOC 40MG 14 F.C.TAB.", "SPASMO-DIGESTIN 30 TABS.")

# Locate every row (product list starts at row 4) whose name matches one of the
# removed products.
$rowsToDelete = New-Object System.Collections.ArrayList
for ($r = 4; $r -le 104; $r++) {
    $v = $ws.Cells($r, 2).Value()
    foreach ($name in $targets) {
        if ($v -eq $name) {
            [void]$rowsToDelete.Add($r)
        }
    }
}

# Delete from the bottom up so the row numbers collected above stay valid.
$sortedRowsDesc = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRowsDesc) {
    $ws.Rows($r).Delete()
}

# Walk the remaining product rows: renumber the sequential index column (A),
# and collect column L ("سعر البيع") to refresh the grand-total cell below the
# table.
$dataStartRow = 4
$n = 0
$prices = New-Object System.Collections.ArrayList
$r = $dataStartRow
while ($true) {
    $name = $ws.Cells($r, 2).Value()
    if (($name -eq $null) -or ($name -eq "")) {
        break
    }
    $n = $n + 1
    $ws.Cells($r, 1).Value = $n
    $priceVal = $ws.Cells($r, 12).Value()
    if ($priceVal -ne $null) {
        [void]$prices.Add($priceVal)
    }
    $r = $r + 1
}
$totalRow = $r

$sortedPrices = $prices | Sort-Object
$total = 0.0
foreach ($v in $sortedPrices) {
    $total = $total + $v
}
$ws.Cells($totalRow, 11).Value = $total
